$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new header columns (AD, AE, AF) matching the existing header
# style used across row 1 (bold, centered, thin border) by copying the
# format from the last existing header cell (AC1).
$ws.Range("AC1").Copy($ws.Range("AD1:AF1"))
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins/Losses/Ties) for every data row.
$lastRow = 55
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 104
    $ws.Cells.Item($r, 31).Value = 58
    $ws.Cells.Item($r, 32).Value = 0
}
